$wb = $excel.ActiveWorkbook

# --- Add the new worksheet "News_Default_Kanal" at the end of the tab
#     strip (after "List_Tab_Menu_Exclusive"), which also makes it the
#     active sheet / activeTab, exactly like the target workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "News_Default_Kanal"

# --- Header row (bold, matches the "Index"/"Tab Name"-style header used
#     by the sibling "List_Tab_Menu_Exclusive" sheet)
$ws.Range("A1").Value = "Index"
$ws.Range("B1").Value = "Value"
$ws.Range("A1:B1").Font.Bold = $true

# --- Data rows. The leading apostrophe forces these numeric-looking
#     strings to be stored as text (quote-prefixed), reusing the shared
#     strings "0"/"1"/"2" already present in the workbook.
$ws.Range("A2").Value = "'0"
$ws.Range("B2").Value = "Berita Utama"

$ws.Range("A3").Value = "'1"
$ws.Range("B3").Value = "Terkini"

$ws.Range("A4").Value = "'2"
$ws.Range("B4").Value = "Populer"

# --- Column B width (best-fit in the source workbook)
$ws.Columns("B").ColumnWidth = 11.7265625

# --- Selection state left behind on the new (now active) sheet
$ws.Range("C1:C1048576").Select() | Out-Null

# --- Page setup used by the new sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
